# Daily attendance processing - 2026-01-04 06:04:42
# Swap the "Recorded By" name order from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" for every row in column G where it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
